$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new trade row (row 9) ---
$ws.Range("A9").Value = 42654.746469907404
# Match the date formatting/style already used by the other Date cells (column A)
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)

$ws.Range("B9").Value = $true
$ws.Range("C9").Value = 9855.2999999999993
$ws.Range("D9").Value = 9840.5400000000009
$ws.Range("E9").Value = 104.43
$ws.Range("F9").Value = 104.74
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = 0.3
$ws.Range("I9").Value = $false

# --- Re-fit column widths to account for the new row's content (mirrors the
# bestFit recalculation Excel performs), translating from the raw OOXML
# column width units to the ColumnWidth property's units. ---
$offset = 0.8333333333333334
$ws.Columns.Item(1).ColumnWidth = 15.375 - $offset
$ws.Columns.Item(2).ColumnWidth = 8.25 - $offset
$ws.Columns.Item(3).ColumnWidth = 7.875 - $offset
$ws.Columns.Item(4).ColumnWidth = 11.25 - $offset
$ws.Columns.Item(5).ColumnWidth = 10.875 - $offset
$ws.Columns.Item(6).ColumnWidth = 7 - $offset
$ws.Columns.Item(7).ColumnWidth = 10.375 - $offset
$ws.Columns.Item(8).ColumnWidth = 14.625 - $offset
$ws.Columns.Item(9).ColumnWidth = 11.875 - $offset
